# "Generate Report for Handback" -- refresh the handback timestamps for the
# ea45270a-3667-44e7-b006-483281ce8c47 row across all report sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# ea45270a-... file, row 4.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-29 10:46:14"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) for the ea45270a-... entry, row 4.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-29 10:46:08"
$wsZhCn.Range("K4").Value = "2016-08-29 10:46:26"

# de-de sheet: same two columns for the ea45270a-... entry, row 4.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-29 10:46:14"
$wsDeDe.Range("K4").Value = "2016-08-29 10:46:33"
